# Update "want-to-go" counts (column F) on the "展览" (id 1) and
# "全部类型" (id 4) worksheets to reflect the new scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览": rows 2-12, column F ---------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    2  = 115
    3  = 7542
    4  = 284
    5  = 16
    6  = 456
    7  = 4157
    8  = 324
    10 = 277
    11 = 664
    12 = 150
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoUpdates[$row]
}

# --- Sheet "全部类型": rows 2-15, column F ------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 115
    4  = 7542
    6  = 284
    7  = 16
    8  = 456
    9  = 4157
    10 = 324
    12 = 277
    13 = 664
    15 = 150
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
